# Apply the Mon Mar 25 18:50:38 UTC 2024 cryptos-list refresh (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values such as "0.999" or "57.70" would otherwise be auto-coerced to numbers
# (losing the trailing zero / becoming a float) by plain Value assignment, so
# for any new value that parses as a number we briefly force Text format, write
# it, then restore the cell's original style so formatting is left untouched.
function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range('D2').Value = '70.689.67'
$ws.Range('E2').Value = '  +7.42%  '

$ws.Range('D3').Value = '3.629.84'
$ws.Range('E3').Value = '  +7.26%  '

$ws.Range('E4').Value = '  +0.01%  '

Set-TextValue 'D5' '594.05'
$ws.Range('E5').Value = '  +5.68%  '

Set-TextValue 'D6' '192.14'
$ws.Range('E6').Value = '  +9.16%  '

$ws.Range('E7').Value = '  +3.10%  '

$ws.Range('D8').Value = '3.606.27'
$ws.Range('E8').Value = '  +6.74%  '

Set-TextValue 'D9' '0.999'
$ws.Range('E9').Value = '  -0.05%  '

$ws.Range('E10').Value = '  +4.16%  '

$ws.Range('E11').Value = '  +4.78%  '

Set-TextValue 'D12' '57.70'
$ws.Range('E12').Value = '  +7.40%  '

Set-TextValue 'D13' '0.0000295'
$ws.Range('E13').Value = '  +6.23%  '

Set-TextValue 'D14' '9.74'
$ws.Range('E14').Value = '  +5.84%  '

$ws.Range('D15').Value = '4.209.54'
$ws.Range('E15').Value = '  +7.47%  '

$ws.Range('D16').Value = '3.628.14'
$ws.Range('E16').Value = '  +7.33%  '

Set-TextValue 'D17' '19.39'
$ws.Range('E17').Value = '  +6.50%  '

$ws.Range('D18').Value = '70.547.86'
$ws.Range('E18').Value = '  +7.42%  '

Set-TextValue 'D19' '12.59'
$ws.Range('E19').Value = '  +6.39%  '

$ws.Range('E20').Value = '  +1.21%  '

$ws.Range('E21').Value = '  +5.50%  '

Set-TextValue 'D22' '492.18'
$ws.Range('E22').Value = '  +5.60%  '

$ws.Range('E23').Value = '  +12.58%  '

Set-TextValue 'D24' '16.73'
$ws.Range('E24').Value = '  +16.81%  '

Set-TextValue 'D25' '4.45'
$ws.Range('E25').Value = '  +8.79%  '

Set-TextValue 'D26' '90.63'
$ws.Range('E26').Value = '  +1.09%  '

$ws.Range('E27').Value = '  +6.52%  '

Set-TextValue 'D28' '11.20'
$ws.Range('E28').Value = '  +5.70%  '

Set-TextValue 'D29' '9.39'
$ws.Range('E29').Value = '  +7.89%  '

Set-TextValue 'D30' '32.38'
$ws.Range('E30').Value = '  +4.28%  '

Set-TextValue 'D31' '7.65'
$ws.Range('E31').Value = '  +16.05%  '

Set-TextValue 'D32' '12.24'
$ws.Range('E32').Value = '  +7.14%  '

Set-TextValue 'D33' '613.85'
$ws.Range('E33').Value = '  +5.98%  '

Set-TextValue 'D34' '65.48'
$ws.Range('E34').Value = '  +5.72%  '

$ws.Range('E35').Value = '  +7.56%  '

$ws.Range('E36').Value = '  +12.14%  '

$ws.Range('E37').Value = '  +4.41%  '

$ws.Range('E38').Value = '  +0.03%  '

Set-TextValue 'D39' '37.93'
$ws.Range('E39').Value = '  +5.39%  '

$ws.Range('E40').Value = '  +7.29%  '

Set-TextValue 'D41' '3.67'
$ws.Range('E41').Value = '  +1.27%  '

$ws.Range('D42').Value = '3.363.86'
$ws.Range('E42').Value = '  +8.73%  '

Set-TextValue 'D43' '3.06'
$ws.Range('E43').Value = '  +8.06%  '

$ws.Range('E44').Value = '  +7.03%  '

Set-TextValue 'D45' '2.65'
$ws.Range('E45').Value = '  +8.71%  '

Set-TextValue 'D46' '3.38'
$ws.Range('E46').Value = '  +7.33%  '

$ws.Range('E47').Value = '  +2.93%  '

$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D48' '9.13'
$ws.Range('E48').Value = '  +7.67%  '

$ws.Range('B49').Value = 'LidoDAOToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D49' '3.36'
$ws.Range('E49').Value = '  +7.37%  '

$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D50' '2.73'
$ws.Range('E50').Value = '  +9.81%  '

Set-TextValue 'D51' '0.999'
$ws.Range('E51').Value = '  +0.15%  '
